$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.512.22"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "'2.298.58"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'537.18"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").Value = "'132.04"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +3.06%  "
$ws.Range("D9").Value = "'2.298.18"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "'0.100"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "'0.333"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "'23.72"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "'2.710.19"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "'58.424.89"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "'2.306.27"
$ws.Range("E18").Value = "  -2.37%  "
$ws.Range("D19").Value = "'10.57"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").Value = "'315.98"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "'6.62"
$ws.Range("E22").Value = "  +2.45%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'63.07"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D27").Value = "'7.93"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("D29").Value = "'171.06"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("D30").Value = "'1.71"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").Value = "'0.0₃0724"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").Value = "'1.09"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").Value = "'5.81"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'17.87"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'1.24"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "'290.05"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("D42").Value = "'140.62"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "'3.45"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").Value = "'0.0950"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'0.0495"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").Value = "'0.555"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "'18.24"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").Value = "'0.0210"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  +0.26%  "
